$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.294
$ws.Range("A3").Value = -21.932
$ws.Range("E3").Value = 16.442
$ws.Range("E12").Value = 17.845
$ws.Range("A14").Value = -21.668
$ws.Range("A21").Value = -19.981
$ws.Range("A23").Value = -19.945
$ws.Range("E24").Value = 16.988
$ws.Range("A25").Value = -21.626
$ws.Range("C25").Value = -11.969
$ws.Range("E25").Value = 17
$ws.Range("A26").Value = -21.006
$ws.Range("C27").Value = -13.309
$ws.Range("A29").Value = -21.284
$ws.Range("C31").Value = -13.163
$ws.Range("C39").Value = -12.858
$ws.Range("C48").Value = -11.754
$ws.Range("E50").Value = 16.331
$ws.Range("C51").Value = -11.133
$ws.Range("C52").Value = -11.601
$ws.Range("A53").Value = -21.888
$ws.Range("E53").Value = 17.036
$ws.Range("C55").Value = -13.175
$ws.Range("C56").Value = -13.537
$ws.Range("A57").Value = -22.07999999999999
$ws.Range("C57").Value = -13.813
$ws.Range("E57").Value = 16.445
$ws.Range("A59").Value = -22.407
$ws.Range("E61").Value = 16.625
$ws.Range("E63").Value = 17.6
$ws.Range("A69").Value = -21.611
$ws.Range("E70").Value = 17.723
$ws.Range("C73").Value = -12.601
$ws.Range("A79").Value = -21.175
$ws.Range("A83").Value = -21.938
$ws.Range("E86").Value = 16.46
$ws.Range("C89").Value = -11.627
$ws.Range("C90").Value = -12.482
$ws.Range("A91").Value = -21.493
$ws.Range("C92").Value = -11.332
$ws.Range("A93").Value = -21.444
$ws.Range("E98").Value = 16.108
$ws.Range("E100").Value = 16.879
$ws.Range("E102").Value = 16.513
